$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 136
$ws.Cells.Item(5, 6).Value = 1993
$ws.Cells.Item(6, 6).Value = 4101
$ws.Cells.Item(7, 6).Value = 534
$ws.Cells.Item(8, 6).Value = 1045
$ws.Cells.Item(9, 6).Value = 657
$ws.Cells.Item(12, 6).Value = 2182
$ws.Cells.Item(14, 6).Value = 654351
$ws.Cells.Item(15, 6).Value = 1613
$ws.Cells.Item(16, 6).Value = 494
$ws.Cells.Item(17, 6).Value = 1441
$ws.Cells.Item(18, 6).Value = 667
$ws.Cells.Item(19, 6).Value = 540
$ws.Cells.Item(20, 6).Value = 1257
$ws.Cells.Item(21, 6).Value = 2198
$ws.Cells.Item(22, 6).Value = 1124
$ws.Cells.Item(23, 6).Value = 2684
$ws.Cells.Item(24, 6).Value = 1542
$ws.Cells.Item(25, 6).Value = 792
$ws.Cells.Item(26, 6).Value = 1523
$ws.Cells.Item(27, 6).Value = 24
$ws.Cells.Item(28, 6).Value = 523
$ws.Cells.Item(29, 6).Value = 1081
$ws.Cells.Item(30, 6).Value = 265
$ws.Cells.Item(33, 6).Value = 76
$ws.Cells.Item(34, 6).Value = 2005
$ws.Cells.Item(35, 6).Value = 1353
$ws.Cells.Item(36, 6).Value = 567
$ws.Cells.Item(37, 6).Value = 1213
$ws.Cells.Item(38, 6).Value = 2302
$ws.Cells.Item(40, 6).Value = 15
$ws.Cells.Item(41, 6).Value = 193
$ws.Cells.Item(42, 6).Value = 2563
$ws.Cells.Item(43, 6).Value = 203
$ws.Cells.Item(44, 6).Value = 979
$ws.Cells.Item(45, 6).Value = 3099
$ws.Cells.Item(46, 6).Value = 339
$ws.Cells.Item(47, 6).Value = 25
$ws.Cells.Item(49, 6).Value = 142

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(9, 6).Value = 100
$ws.Cells.Item(10, 6).Value = 475
$ws.Cells.Item(11, 6).Value = 144491
$ws.Cells.Item(12, 6).Value = 144491
$ws.Cells.Item(17, 6).Value = 96
$ws.Cells.Item(18, 6).Value = 226
$ws.Cells.Item(19, 6).Value = 332
$ws.Cells.Item(21, 6).Value = 409
$ws.Cells.Item(22, 6).Value = 409
$ws.Cells.Item(23, 6).Value = 119
$ws.Cells.Item(24, 6).Value = 80
$ws.Cells.Item(25, 6).Value = 96
$ws.Cells.Item(26, 6).Value = 89
$ws.Cells.Item(27, 6).Value = 537
$ws.Cells.Item(31, 6).Value = 58
$ws.Cells.Item(32, 6).Value = 334
$ws.Cells.Item(35, 6).Value = 49
$ws.Cells.Item(38, 6).Value = 111
$ws.Cells.Item(39, 6).Value = 10

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 3120
$ws.Cells.Item(5, 6).Value = 240
$ws.Cells.Item(7, 6).Value = 820
$ws.Cells.Item(8, 6).Value = 1163
$ws.Cells.Item(9, 6).Value = 634
$ws.Cells.Item(10, 6).Value = 1585
$ws.Cells.Item(12, 6).Value = 73
$ws.Cells.Item(13, 6).Value = 1852

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 820
$ws.Cells.Item(3, 6).Value = 1163
$ws.Cells.Item(4, 6).Value = 634
$ws.Cells.Item(6, 6).Value = 1585
$ws.Cells.Item(8, 6).Value = 136
$ws.Cells.Item(9, 6).Value = 1993
$ws.Cells.Item(10, 6).Value = 73
$ws.Cells.Item(11, 6).Value = 1852
$ws.Cells.Item(12, 6).Value = 4102
$ws.Cells.Item(13, 6).Value = 534
$ws.Cells.Item(14, 6).Value = 657
$ws.Cells.Item(16, 6).Value = 2182
$ws.Cells.Item(18, 6).Value = 654361
$ws.Cells.Item(19, 6).Value = 100
$ws.Cells.Item(20, 6).Value = 475
$ws.Cells.Item(21, 6).Value = 1613
$ws.Cells.Item(22, 6).Value = 144491
$ws.Cells.Item(23, 6).Value = 1441
$ws.Cells.Item(24, 6).Value = 667
$ws.Cells.Item(25, 6).Value = 540
$ws.Cells.Item(26, 6).Value = 1257
$ws.Cells.Item(27, 6).Value = 2198
$ws.Cells.Item(28, 6).Value = 1124
$ws.Cells.Item(29, 6).Value = 2684
$ws.Cells.Item(30, 6).Value = 1542
$ws.Cells.Item(31, 6).Value = 792
$ws.Cells.Item(33, 6).Value = 1523
$ws.Cells.Item(34, 6).Value = 409
$ws.Cells.Item(35, 6).Value = 523
$ws.Cells.Item(36, 6).Value = 119
$ws.Cells.Item(37, 6).Value = 1081
$ws.Cells.Item(39, 6).Value = 76
$ws.Cells.Item(40, 6).Value = 2005
$ws.Cells.Item(41, 6).Value = 1353
$ws.Cells.Item(42, 6).Value = 1213
$ws.Cells.Item(43, 6).Value = 2302
$ws.Cells.Item(45, 6).Value = 334
$ws.Cells.Item(46, 6).Value = 334
$ws.Cells.Item(48, 6).Value = 2563
$ws.Cells.Item(49, 6).Value = 203
$ws.Cells.Item(50, 6).Value = 979
$ws.Cells.Item(51, 6).Value = 3099
$ws.Cells.Item(52, 6).Value = 142
